$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are treated as text so strings like "209.93"
# or "1.00" are not silently re-interpreted as numbers (losing formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '89.000.60'
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').Value = '3.040.20'
$ws.Range('E3').Value = '  -1.60%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').Value = '209.93'
$ws.Range('E5').Value = '  -0.82%  '
$ws.Range('D6').Value = '610.34'
$ws.Range('E6').Value = '  -3.28%  '
$ws.Range('D7').Value = '0.361'
$ws.Range('E7').Value = '  -7.09%  '
$ws.Range('D8').Value = '0.876'
$ws.Range('E8').Value = '  +21.77%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '3.035.30'
$ws.Range('E10').Value = '  -1.66%  '
$ws.Range('D11').Value = '0.659'
$ws.Range('E11').Value = '  +19.31%  '
$ws.Range('E12').Value = '  +3.55%  '
$ws.Range('D13').Value = '0.0000236'
$ws.Range('E13').Value = '  -4.81%  '
$ws.Range('E14').Value = '  +1.74%  '
$ws.Range('D15').Value = '88.442.79'
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').Value = '31.66'
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '3.595.45'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').Value = '3.064.15'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').Value = '3.37'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').Value = '0.0000207'
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').Value = '13.34'
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('D22').Value = '423.25'
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').Value = '4.97'
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('D24').Value = '8.05'
$ws.Range('E24').Value = '  -2.85%  '
$ws.Range('D25').Value = '5.40'
$ws.Range('E25').Value = '  +3.13%  '
$ws.Range('D26').Value = '83.07'
$ws.Range('E26').Value = '  +5.35%  '
$ws.Range('D27').Value = '11.60'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').Value = '3.194.33'
$ws.Range('E28').Value = '  -1.74%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').Value = '1.06'
$ws.Range('E30').Value = '  +95.11%  '
$ws.Range('D31').Value = '0.162'
$ws.Range('E31').Value = '  +3.65%  '
$ws.Range('D32').Value = '8.16'
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('D33').Value = '499.56'
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('D34').Value = '3.56'
$ws.Range('E34').Value = '  -9.78%  '
$ws.Range('D35').Value = '6.58'
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('D36').Value = '1.79'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').Value = '22.50'
$ws.Range('E37').Value = '  +4.16%  '
$ws.Range('D38').Value = '1.23'
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').Value = '22.21'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E40').Value = '  +4.99%  '
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('E43').Value = '  +11.27%  '
$ws.Range('D44').Value = '0.360'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').Value = '1.81'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').Value = '146.23'
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('D47').Value = '43.30'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '0.0675'
$ws.Range('E48').Value = '  +11.83%  '
$ws.Range('D49').Value = '4.04'
$ws.Range('E49').Value = '  +3.74%  '
$ws.Range('D50').Value = '1.20'
$ws.Range('E50').Value = '  +3.12%  '
$ws.Range('D51').Value = '155.56'
$ws.Range('E51').Value = '  -5.10%  '
